$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 457, pushing existing rows 457:555 down to 459:557
$ws.Rows.Item(457).Resize(2).Insert()

# Populate the first new row (457)
$ws.Cells.Item(457,1).Value = 4
$ws.Cells.Item(457,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(457,3).Value = "Los Lagos"
$ws.Cells.Item(457,4).Value = 45244
$ws.Cells.Item(457,5).Value = 10
$ws.Cells.Item(457,6).Value = 100112003
$ws.Cells.Item(457,7).Value = "Ajo"
$ws.Cells.Item(457,8).Value = "Chino"
$ws.Cells.Item(457,9).Value = "Primera"
$ws.Cells.Item(457,10).Value = 240
$ws.Cells.Item(457,11).Value = 26000
$ws.Cells.Item(457,12).Value = 26000
$ws.Cells.Item(457,13).Value = 26000
$ws.Cells.Item(457,14).Value = "$/caja 10 kilos"
$ws.Cells.Item(457,15).Value = "China"
$ws.Cells.Item(457,16).Value = 2600
$ws.Cells.Item(457,17).Value = 10
$ws.Cells.Item(457,18).Value = "Hortaliza"

# Populate the second new row (458)
$ws.Cells.Item(458,1).Value = 4
$ws.Cells.Item(458,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(458,3).Value = "Los Lagos"
$ws.Cells.Item(458,4).Value = 45244
$ws.Cells.Item(458,5).Value = 10
$ws.Cells.Item(458,6).Value = 100112003
$ws.Cells.Item(458,7).Value = "Ajo"
$ws.Cells.Item(458,8).Value = "Chino"
$ws.Cells.Item(458,9).Value = "Primera"
$ws.Cells.Item(458,10).Value = 90
$ws.Cells.Item(458,11).Value = 27000
$ws.Cells.Item(458,12).Value = 27000
$ws.Cells.Item(458,13).Value = 27000
$ws.Cells.Item(458,14).Value = "$/malla 10 kilos"
$ws.Cells.Item(458,15).Value = "China"
$ws.Cells.Item(458,16).Value = 2700
$ws.Cells.Item(458,17).Value = 10
$ws.Cells.Item(458,18).Value = "Hortaliza"
